$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.041.31"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.640.92"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = $ws.Range("C4").Style
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'215.24"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "'0.5047"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.2575"
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.06432"
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "'19.45"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "'0.07718"
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "1.642.53"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "'4.246"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "1.867.01"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "'0.5454"
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "0.0₅7903"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "'63.58"
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "26.023.79"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'203.89"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").Value = "'4.293"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").Value = "'9.995"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").Value = "'5.970"
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "'1.009"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'1.932"
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = "  +9.48%  "
$ws.Range("D26").Value = "'141.35"
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").Value = "'0.1153"
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "'6.737"
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").Value = "'0.05057"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").Value = "'3.253"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("D33").Value = "'3.191"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").Value = "'2.340"
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").Value = "'0.8950"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("E37").Value = "  -5.41%  "
$ws.Range("D38").Value = "'0.5630"
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "1.148.56"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").Value = "'2.564"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'5.676"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'0.8113"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").Value = "'99.85"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "1.777.69"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").Value = "'0.4529"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'54.92"
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  -0.97%  "
